$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = "Espesor e m"
    4  = "Responsable"
    5  = "teléfono"
    6  = "profesión"
    7  = "Velocidad del viento ms"
    8  = "nif"
    9  = "dia/s de visita"
    10 = "Longitud invernadero m-0"
    11 = "K kWm2 C-0"
    12 = "nombre 1"
    13 = "Conductividad térmica A Wm C-0"
    14 = "a"
    15 = "Firma1"
    16 = "referencia catastral"
    18 = "mes"
    19 = "NIF"
    20 = "t C-2"
    21 = "dirección postal"
    22 = "t C-1"
    23 = "año"
    24 = "K kWm2 C"
    25 = "Conductividad térmica A Wm C"
    26 = "ciudad"
    27 = "Espesor e m-0"
    28 = "Velocidad del viento ms-0"
    29 = "mail"
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 1).Value = $values[$row]
}
